$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "@1 Ohm"
$ws.Range("B16").Value = 122

[void]$ws.Range("B17").Select()
